$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '26.996.68'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -1.63%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.821.64'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -1.14%  '

$ws.Range('E4').Value = '  -0.54%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '309.54'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.84%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.008'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.52%  '

$ws.Range('E7').Value = '  -2.97%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3639'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -1.95%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07283'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -2.49%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.8644'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -2.74%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '19.82'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -3.36%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.885.48'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +2.24%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.07599'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +2.86%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '93.21'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -0.36%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '5.323'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -2.88%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '6.492'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -1.62%  '

$ws.Range('E17').Value = '  -0.65%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.000008624'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -2.55%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '27.380.14'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -0.30%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '14.47'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -2.54%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.158'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -3.62%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '10.58'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -1.54%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.115.02'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +2.04%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '151.62'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -0.53%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.853'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -2.59%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '18.20'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -2.42%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.087'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -3.92%  '

$ws.Range('E29').Value = '  -3.62%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '115.84'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -1.89%  '

$ws.Range('E31').Value = '  -1.05%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '2.951'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.08%  '

$ws.Range('B33').Value = 'ARBITRUM'
$ws.Range('C33').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.140'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -3.43%  '

$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.7269'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -4.33%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.424'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -3.19%  '

$ws.Range('E36').Value = '  -0.49%  '

$ws.Range('E37').Value = '  +5.11%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.05275'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -1.63%  '

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.074'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -2.90%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.01917'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -2.71%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.926'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -2.43%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '7.164'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -2.07%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.5200'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -3.12%  '

$ws.Range('E44').Value = '  -2.03%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '8.259'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -3.57%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.4856'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -2.50%  '

$ws.Range('E47').Value = '  -0.55%  '

$ws.Range('E48').Value = '  -5.21%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '103.27'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -1.95%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.630'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -3.36%  '

$ws.Range('E51').Value = '  -1.60%  '
